$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024 April 06 4:16:05 AM", "ASP_SERVER", "Packet sent. Type: ReadyImage"),
    @("2024 April 06 4:16:10 AM", "TCP_SERVER", "Images received. Count: 6"),
    @("2024 April 06 4:16:10 AM", "ASP_SERVER", "Packet sent. Type: ReadyPost"),
    @("2024 April 06 4:16:10 AM", "TCP_SERVER", "Posts received. Count: 4")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
